$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "61.629.97"
$ws.Range("E2").Value = "  -5.90%  "
$ws.Range("D3").Value = "3.129.59"
$ws.Range("E3").Value = "  -7.61%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue "D5" "509.25"
$ws.Range("E5").Value = "  -3.62%  "
Set-TextValue "D6" "167.23"
$ws.Range("E6").Value = "  -10.52%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "3.126.06"
$ws.Range("E9").Value = "  -7.60%  "
Set-TextValue "D10" "0.583"
$ws.Range("E10").Value = "  -6.95%  "
Set-TextValue "D11" "51.12"
$ws.Range("E11").Value = "  -12.95%  "
$ws.Range("E12").Value = "  -5.64%  "
Set-TextValue "D13" "0.0000244"
$ws.Range("E13").Value = "  -4.40%  "
Set-TextValue "D14" "8.63"
$ws.Range("E14").Value = "  -6.83%  "
$ws.Range("D15").Value = "3.628.84"
$ws.Range("E15").Value = "  -7.41%  "
$ws.Range("E16").Value = "  -8.75%  "
$ws.Range("D17").Value = "3.127.64"
$ws.Range("E17").Value = "  -7.57%  "
$ws.Range("D18").Value = "61.496.79"
$ws.Range("E18").Value = "  -5.59%  "
Set-TextValue "D19" "16.68"
$ws.Range("E19").Value = "  -4.60%  "
$ws.Range("E20").Value = "  -5.49%  "
Set-TextValue "D21" "0.932"
$ws.Range("E21").Value = "  -4.36%  "
Set-TextValue "D22" "355.42"
$ws.Range("E22").Value = "  -4.76%  "
Set-TextValue "D23" "78.88"
$ws.Range("E23").Value = "  -3.66%  "
$ws.Range("E24").Value = "  -3.29%  "
Set-TextValue "D25" "10.61"
$ws.Range("E25").Value = "  -2.36%  "
Set-TextValue "D26" "6.08"
$ws.Range("E26").Value = "  +3.94%  "
Set-TextValue "D27" "3.78"
$ws.Range("E27").Value = "  +2.24%  "
$ws.Range("E28").Value = "  -4.95%  "
Set-TextValue "D29" "10.83"
$ws.Range("E29").Value = "  -6.11%  "
Set-TextValue "D30" "7.89"
$ws.Range("E30").Value = "  -7.62%  "
Set-TextValue "D31" "633.79"
$ws.Range("E31").Value = "  -6.54%  "
Set-TextValue "D32" "27.49"
$ws.Range("E32").Value = "  -7.39%  "
Set-TextValue "D33" "6.19"
$ws.Range("E33").Value = "  -8.04%  "
Set-TextValue "D34" "10.92"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -5.44%  "
Set-TextValue "D37" "55.29"
$ws.Range("E37").Value = "  -9.69%  "
Set-TextValue "D38" "35.11"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("E39").Value = "  -5.83%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("D41").Value = "0.0₃0675"
$ws.Range("E41").Value = "  +8.90%  "
Set-TextValue "D42" "0.119"
$ws.Range("E42").Value = "  -6.70%  "
$ws.Range("D43").Value = "2.782.98"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  +2.28%  "
Set-TextValue "D45" "2.59"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D48" "2.90"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D49" "2.46"
$ws.Range("E49").Value = "  -11.49%  "
Set-TextValue "D50" "131.37"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("E51").Value = "  -5.08%  "
